$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 / C3: top_four_third_teams
$ws.Range("C2").Value = "['Czech Republic', 'Northern Ireland', 'Portugal', 'Slovakia']"
$ws.Range("C3").Value = "['Czech Republic', 'Northern Ireland', 'Portugal', 'Slovakia']"

# D2 / D3: last_two_third_teams
$ws.Range("D2").Value = "['Romania', 'Sweden']"
$ws.Range("D3").Value = "['Romania', 'Sweden']"

# D7 / D8 / D9: last_two_third_teams - reorder
$ws.Range("D7").Value = "['Portugal', 'Sweden']"
$ws.Range("D8").Value = "['Portugal', 'Sweden']"
$ws.Range("D9").Value = "['Portugal', 'Sweden']"

# C30: top_four_third_teams
$ws.Range("C30").Value = "['Austria', 'Croatia', 'Finland', 'Hungary']"
# D30: last_two_third_teams
$ws.Range("D30").Value = "['Switzerland', 'Spain']"

# D36 / D37: reorder
$ws.Range("D36").Value = "['Denmark', 'Slovakia']"
$ws.Range("D37").Value = "['Denmark', 'Slovakia']"

# C38
$ws.Range("C38").Value = "['Denmark', 'Hungary', 'Switzerland', 'Ukraine']"
# D38
$ws.Range("D38").Value = "['Slovakia', 'Croatia']"

# C39
$ws.Range("C39").Value = "['Denmark', 'Germany', 'Switzerland', 'Ukraine']"
# D39
$ws.Range("D39").Value = "['Slovakia', 'Croatia']"

# F40: change_flag
$ws.Range("F40").Value = 0

# G40..G59: change_count decremented by 1
$ws.Range("G40").Value = 7
$ws.Range("G41").Value = 7
$ws.Range("G42").Value = 7
$ws.Range("G43").Value = 7
$ws.Range("G44").Value = 7
$ws.Range("G45").Value = 8
$ws.Range("G46").Value = 8
$ws.Range("G47").Value = 9
$ws.Range("G48").Value = 10
$ws.Range("G49").Value = 10
$ws.Range("G50").Value = 11
$ws.Range("G51").Value = 12
$ws.Range("G52").Value = 12
$ws.Range("G53").Value = 12
$ws.Range("G54").Value = 13
$ws.Range("G55").Value = 13
$ws.Range("G56").Value = 13
$ws.Range("G57").Value = 13
$ws.Range("G58").Value = 14
$ws.Range("G59").Value = 14
